# Scheduled-runner style data refresh: updates market-price columns
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ] -> H:N)
# for a handful of leve rows across several job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3250.4119
$ws.Range("I86").Value = 4072.8572
$ws.Range("K86").Value = 4072.8572
$ws.Range("M86").Value = -2949.8572

$ws.Range("H89").Value = 3250.4119
$ws.Range("I89").Value = 4072.8572
$ws.Range("K89").Value = 20364.286
$ws.Range("M89").Value = -14748.286

$ws.Range("H112").Value = 31084.234
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 31980.727
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 95942.181
$ws.Range("M112").Value = -3392
$ws.Range("N112").Value = -98158.181

$ws.Range("H137").Value = 2150.1025
$ws.Range("I137").Value = 2019.4736
$ws.Range("J137").Value = 2504.6667
$ws.Range("K137").Value = 6058.4208
$ws.Range("L137").Value = 7514.000100000001
$ws.Range("M137").Value = -3508.4208
$ws.Range("N137").Value = -12614.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 948
$ws.Range("I61").Value = 909.03705
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 909.03705
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -697.03705
$ws.Range("N61").Value = -2424

$ws.Range("H74").Value = 1798.2059
$ws.Range("I74").Value = 764.13635
$ws.Range("J74").Value = 3694
$ws.Range("K74").Value = 764.13635
$ws.Range("L74").Value = 3694
$ws.Range("M74").Value = 109.86365
$ws.Range("N74").Value = -5442

$ws.Range("H77").Value = 1798.2059
$ws.Range("I77").Value = 764.13635
$ws.Range("J77").Value = 3694
$ws.Range("K77").Value = 3820.68175
$ws.Range("L77").Value = 18470
$ws.Range("M77").Value = 547.3182500000003
$ws.Range("N77").Value = -27206

$ws.Range("H136").Value = 948
$ws.Range("I136").Value = 909.03705
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 2727.11115
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -177.1111500000002
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2923.439
$ws.Range("I86").Value = 2769.8518
$ws.Range("J86").Value = 3219.6428
$ws.Range("K86").Value = 2769.8518
$ws.Range("L86").Value = 3219.6428
$ws.Range("M86").Value = -1646.8518
$ws.Range("N86").Value = -5465.6428

$ws.Range("H89").Value = 2923.439
$ws.Range("I89").Value = 2769.8518
$ws.Range("J89").Value = 3219.6428
$ws.Range("K89").Value = 13849.259
$ws.Range("L89").Value = 16098.214
$ws.Range("M89").Value = -8233.259
$ws.Range("N89").Value = -27330.214

$ws.Range("H134").Value = 1126.0952
$ws.Range("I134").Value = 1081.7
$ws.Range("J134").Value = 2014
$ws.Range("K134").Value = 3245.1
$ws.Range("L134").Value = 6042
$ws.Range("M134").Value = -710.1000000000004
$ws.Range("N134").Value = -11112

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22354.787
$ws.Range("I31").Value = 25804.643
$ws.Range("J31").Value = 14728.789
$ws.Range("K31").Value = 25804.643
$ws.Range("L31").Value = 14728.789
$ws.Range("M31").Value = -25509.643
$ws.Range("N31").Value = -15318.789

$ws.Range("H34").Value = 22354.787
$ws.Range("I34").Value = 25804.643
$ws.Range("J34").Value = 14728.789
$ws.Range("K34").Value = 25804.643
$ws.Range("L34").Value = 14728.789
$ws.Range("M34").Value = -25602.643
$ws.Range("N34").Value = -15132.789

$ws.Range("H58").Value = 1192.8529
$ws.Range("I58").Value = 1183.1212
$ws.Range("K58").Value = 1183.1212
$ws.Range("M58").Value = -980.1212

$ws.Range("H99").Value = 1882.2222
$ws.Range("I99").Value = 1537.3334
$ws.Range("J99").Value = 2054.6667
$ws.Range("K99").Value = 1537.3334
$ws.Range("L99").Value = 2054.6667
$ws.Range("M99").Value = -39.33339999999998
$ws.Range("N99").Value = -5050.6667

$ws.Range("H126").Value = 1882.2222
$ws.Range("I126").Value = 1537.3334
$ws.Range("J126").Value = 2054.6667
$ws.Range("K126").Value = 4612.0002
$ws.Range("L126").Value = 6164.000100000001
$ws.Range("M126").Value = -2142.0002
$ws.Range("N126").Value = -11104.0001

$ws.Range("H132").Value = 897.717
$ws.Range("I132").Value = 642.7805
$ws.Range("J132").Value = 1768.75
$ws.Range("K132").Value = 1928.3415
$ws.Range("L132").Value = 5306.25
$ws.Range("M132").Value = 601.6585
$ws.Range("N132").Value = -10366.25

$ws.Range("H134").Value = 1674.7142
$ws.Range("I134").Value = 1256.3914
$ws.Range("J134").Value = 2476.5
$ws.Range("K134").Value = 3769.1742
$ws.Range("L134").Value = 7429.5
$ws.Range("M134").Value = -1234.1742
$ws.Range("N134").Value = -12499.5

$ws.Range("H136").Value = 1192.8529
$ws.Range("I136").Value = 1183.1212
$ws.Range("K136").Value = 3549.3636
$ws.Range("M136").Value = -999.3636000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 110
$ws.Range("I23").Value = 19.333334
$ws.Range("J23").Value = 170.44444
$ws.Range("K23").Value = 58.000002
$ws.Range("L23").Value = 511.33332
$ws.Range("M23").Value = 176.999998
$ws.Range("N23").Value = -981.33332

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2403.5
$ws.Range("I102").Value = 2123.2222
$ws.Range("J102").Value = 2908
$ws.Range("K102").Value = 2123.2222
$ws.Range("L102").Value = 2908
$ws.Range("M102").Value = -501.2222000000002
$ws.Range("N102").Value = -6152

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1850.7391
$ws.Range("I132").Value = 1575.8125
$ws.Range("J132").Value = 2479.1428
$ws.Range("K132").Value = 4727.4375
$ws.Range("L132").Value = 7437.428400000001
$ws.Range("M132").Value = -2197.4375
$ws.Range("N132").Value = -12497.4284

$ws.Range("H136").Value = 1936.8904
$ws.Range("I136").Value = 1340.1538
$ws.Range("J136").Value = 3414.524
$ws.Range("K136").Value = 4020.4614
$ws.Range("L136").Value = 10243.572
$ws.Range("M136").Value = -1470.4614
$ws.Range("N136").Value = -15343.572

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I132").Value = 1008.0741
$ws.Range("J132").Value = 1706.1177
$ws.Range("K132").Value = 3024.2223
$ws.Range("L132").Value = 5118.3531
$ws.Range("M132").Value = -494.2223000000004
$ws.Range("N132").Value = -10178.3531

$ws.Range("H136").Value = 681.14813
$ws.Range("I136").Value = 385.05
$ws.Range("J136").Value = 1527.1428
$ws.Range("K136").Value = 1155.15
$ws.Range("L136").Value = 4581.428400000001
$ws.Range("M136").Value = 1394.85
$ws.Range("N136").Value = -9681.428400000001
